$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7218
$ws1.Range("F3").Value = 405
$ws1.Range("F4").Value = 120
$ws1.Range("F5").Value = 171
$ws1.Range("F7").Value = 92
$ws1.Range("F8").Value = 605

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7218
$ws4.Range("F3").Value = 405
$ws4.Range("F4").Value = 7
$ws4.Range("F5").Value = 120
$ws4.Range("F6").Value = 171
$ws4.Range("F9").Value = 92
$ws4.Range("F10").Value = 605
